# Wossal-web "template-admin" CSV import template refresh.
# - Rename sheet
# - Re-label / reorder header row, refresh sample row values
# - Re-theme (Google-ish palette, Arial everywhere, Poppins highlight cell)
# - Re-size columns, drop unused trailing rows' placeholder formatting

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Sheet title
# ---------------------------------------------------------------------
$ws.Name = "Template fichier CSV pour impor"

# ---------------------------------------------------------------------
# Header row (row 1) -- new wording / new trailing "Catégorie
# socioprofessionnelle" column
# ---------------------------------------------------------------------
$ws.Range("A1").Value = "Nom"
$ws.Range("B1").Value = "Prénom"
$ws.Range("C1").Value = "Adresse mail"
$ws.Range("D1").Value = "Date de naissance"
$ws.Range("E1").Value = "Numéro téléphone"
$ws.Range("F1").Value = "Adresse postale"
$ws.Range("G1").Value = "Fonction"
$ws.Range("H1").Value = "Matricule"
$ws.Range("I1").Value = "Salaire"
$ws.Range("J1").Value = "Catégorie socioprofessionnelle"

# ---------------------------------------------------------------------
# Sample data row (row 2)
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "DIOP"
$ws.Range("B2").Value = "Moussa"
$ws.Range("C2").Value = "papemactarndiayepro+testpro@gmail.com"
$ws.Range("D2").Value = 35417
$ws.Range("D2").NumberFormat = "d/m/yyyy"
$ws.Range("E2").Formula = "=221784563231"
$ws.Range("F2").Value = 9029092
$ws.Range("G2").Value = "Testeur"
$ws.Range("H2").Value = "JKSJK892892"
# I2 used to be a text-formatted formula cell ("=1230303"); drop the old
# Text number format/formula entirely so the new sample value is stored
# as a real number.
$ws.Range("I2").Clear()
$ws.Range("I2").Value = 82828282
$ws.Range("J2").Value = "tests Daly 2003"

# ---------------------------------------------------------------------
# Drop the old placeholder formatting that used to sit under the date
# column / the "Identifiant unique" column on rows 3-6
# ---------------------------------------------------------------------
$ws.Range("D3:D6").Clear()
$ws.Range("J3:J4").Clear()

# ---------------------------------------------------------------------
# Fonts: Calibri -> Arial everywhere
# ---------------------------------------------------------------------
$ws.Cells.Font.Name = "Arial"

# Highlight the new trailing column's sample cell (Poppins, navy text,
# white fill)
$hl = $ws.Range("J2")
$hl.Font.Name = "Poppins"
$hl.Font.Color = 6036998
$hl.Interior.Color = 16777215

# ---------------------------------------------------------------------
# Column widths
# ---------------------------------------------------------------------
$ws.Range("A:C").ColumnWidth = 12.63
$ws.Range("D:D").ColumnWidth = 25.5
$ws.Range("E:E").ColumnWidth = 33.0
$ws.Range("F:F").ColumnWidth = 15.13
$ws.Range("G:G").ColumnWidth = 24.5

# ---------------------------------------------------------------------
# Theme re-colour (Google-ish palette) + Arial theme fonts
# ---------------------------------------------------------------------
$theme = $wb.Theme
$scheme = $theme.ThemeColorScheme
$scheme.Colors(5).RGB = 16024898   # accent1 4285F4
$scheme.Colors(6).RGB = 3490794    # accent2 EA4335
$scheme.Colors(7).RGB = 310523     # accent3 FBBC04
$scheme.Colors(8).RGB = 5482548    # accent4 34A853
$scheme.Colors(9).RGB = 93695      # accent5 FF6D01
$scheme.Colors(10).RGB = 13024582  # accent6 46BDC6
$scheme.Colors(11).RGB = 13391121  # hlink 1155CC
$scheme.Colors(12).RGB = 13391121  # folHlink 1155CC

$fontScheme = $theme.ThemeFontScheme
$fontScheme.MajorFont(1).Name = "Arial"
$fontScheme.MinorFont(1).Name = "Arial"

Write-Host "template-admin refresh applied"
